$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 ("méthode permettant de remettre l'image dans la bonne zone") entirely;
# all following rows shift up by one.
$ws.Rows(6).Delete()

# Fill in previously empty "Temps réel" (D) cells for the first three tasks.
$ws.Range("D2").Value2 = "1H"
$ws.Range("D3").Value2 = "2H"
$ws.Range("D4").Value2 = "1H"

# Fix the typo "Bene" -> "Bénédicte" (row that now sits at 7 after the deletion).
$ws.Range("C7").Value2 = "Bénédicte"

# Fill in the previously empty "Temps provisoire" / "Temps réel" / "Temps réel"(date)
# for the "Tracking" row (now row 9). Borrow the number/cell formatting from
# neighbouring cells that already carry the right style before writing values.
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("E10").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("B9").Value2 = "20mn"
$ws.Range("D9").Value2 = "20mn"
$ws.Range("E9").Value2 = 42060

# Fill in previously empty "Effectué par" / "Temps réel" for the last two rows.
$ws.Range("C13").Value2 = "Bénédicte"
$ws.Range("D13").Value2 = "1H30"
$ws.Range("C14").Value2 = "Bénédicte"
$ws.Range("D14").Value2 = "20mn"

# Update the saved selection to match the author's final cursor position.
$ws.Range("C16").Select()
